$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.035.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.562.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.10'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0855'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.785.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.562.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.036.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("E22").Value = '  +1.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("E30").Value = '  +1.05%  '
$ws.Range("E31").Value = '  +2.83%  '
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.425.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.30%  '
$ws.Range("E36").Value = '  +7.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.32%  '
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.534'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.35%  '
$ws.Range("E40").Value = '  +1.80%  '
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.70'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.699.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("E49").Value = '  +4.99%  '
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("E51").Value = '  +0.29%  '
